$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 262, shifting existing rows 262:267 down to 263:268
$ws.Rows.Item(262).Insert()

# Populate the newly inserted row 262 with the new data point
$ws.Cells.Item(262, 1).Value = 3
$ws.Cells.Item(262, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(262, 3).Value = "Coquimbo"
$ws.Cells.Item(262, 4).Value = 45239
$ws.Cells.Item(262, 4).NumberFormat = $ws.Cells.Item(263, 4).NumberFormat
$ws.Cells.Item(262, 5).Value = 5
$ws.Cells.Item(262, 6).Value = 100112052
$ws.Cells.Item(262, 7).Value = "Albahaca"
$ws.Cells.Item(262, 8).Value = "Sin especificar"
$ws.Cells.Item(262, 9).Value = "Primera"
$ws.Cells.Item(262, 10).Value = 50
$ws.Cells.Item(262, 11).Value = 5000
$ws.Cells.Item(262, 12).Value = 5000
$ws.Cells.Item(262, 13).Value = 5000
$ws.Cells.Item(262, 14).Value = "`$/docena de matas"
$ws.Cells.Item(262, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(262, 16).Value = 833
$ws.Cells.Item(262, 17).Value = 6
$ws.Cells.Item(262, 18).Value = "Hortaliza"
